# Update countries & provincias Spain
# Applies the data refresh (COVID-19 country stats) described by the diff:
#  - Update the "last updated" timestamp in A1
#  - Update numeric stats for several countries (Estados Unidos, Alemania,
#    Suiza, India, Republica de Chipre)
#  - Re-order Eslovaquia so it now appears right before Marruecos (with
#    refreshed numbers), shifting Marruecos down to the old Uruguay row and
#    Uruguay down to the old Eslovaquia row (their own stats unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: last updated timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 10:42"

# --- Row 6: Estados Unidos ---
$ws.Range("B6").Value = 68573
$ws.Range("C6").Value = 362
$ws.Range("D6").Value = 428
$ws.Range("E6").Value = 67109
$ws.Range("G6").Value = 9
$ws.Range("H6").Value = 1036

# --- Row 8: Alemania ---
$ws.Range("B8").Value = 39312
$ws.Range("C8").Value = 1989
$ws.Range("E8").Value = 35543
$ws.Range("G8").Value = 16
$ws.Range("H8").Value = 222

# --- Row 11: Suiza ---
$ws.Range("B11").Value = 11027
$ws.Range("C11").Value = 130
$ws.Range("E11").Value = 10742
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 154

# --- Row 46: India ---
$ws.Range("B46").Value = 693
$ws.Range("C46").Value = 36
$ws.Range("E46").Value = 637

# --- Rows 73-75: Eslovaquia moves ahead of Marruecos ---
# Row 73 now holds Eslovaquia's refreshed data (previously Marruecos' row)
$ws.Range("A73").Value = "Eslovaquia"
$ws.Range("B73").Value = 226
$ws.Range("C73").Value = 10
$ws.Range("D73").Value = 7
$ws.Range("E73").Value = 219
$ws.Range("F73").Value = 2
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0

# Row 74 now holds Marruecos (same stats it had before, just shifted down)
$ws.Range("A74").Value = "Marruecos"
$ws.Range("B74").Value = 225
$ws.Range("C74").Value = 0
$ws.Range("D74").Value = 7
$ws.Range("E74").Value = 212
$ws.Range("F74").Value = 1
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 6

# Row 75 now holds Uruguay (same stats it had before, just shifted down)
$ws.Range("A75").Value = "Uruguay"
$ws.Range("B75").Value = 217
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 217
$ws.Range("F75").Value = 3
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0

# --- Row 90: Republica de Chipre ---
$ws.Range("D90").Value = 4
$ws.Range("E90").Value = 125
